# Update Price (D) and Volume 1h (E) columns for the crypto price table
# A leading apostrophe is used on numeric-looking Price values so Excel
# keeps them as literal text (matching the original inline-string cells)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.496.36'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.917.78'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("D4").Value = '''1.014'
$ws.Range("E4").Value = '  +0.83%  '
$ws.Range("D5").Value = '''325.09'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '''1.012'
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").Value = '''0.4806'
$ws.Range("E7").Value = '  -0.66%  '
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").Value = '''0.08207'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("D11").Value = '''23.45'
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = '1.919.39'
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("D13").Value = '''6.049'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '''7.226'
$ws.Range("E14").Value = '  +1.81%  '
$ws.Range("D15").Value = '''91.39'
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("D16").Value = '''0.06861'
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '''0.00001037'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = '''17.53'
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("E20").Value = '  +0.68%  '
$ws.Range("D21").Value = '29.502.29'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = '''5.664'
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").Value = '''11.87'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '''2.193'
$ws.Range("E24").Value = '  +1.36%  '
$ws.Range("D25").Value = '2.153.36'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").Value = '''156.15'
$ws.Range("D27").Value = '''6.466'
$ws.Range("E27").Value = '  +2.77%  '
$ws.Range("D28").Value = '''20.00'
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").Value = '''2.093'
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").Value = '''120.38'
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("D33").Value = '''5.608'
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("D34").Value = '''3.558'
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("E35").Value = '  -1.86%  '
$ws.Range("D36").Value = '''0.06311'
$ws.Range("E36").Value = '  +3.26%  '
$ws.Range("D37").Value = '''0.02279'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").Value = '''1.179'
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").Value = '''0.5925'
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("D40").Value = '''10.70'
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("D41").Value = '''7.890'
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("D42").Value = '''0.1844'
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").Value = '''2.409'
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("D44").Value = '''1.281'
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("D45").Value = '''12.36'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '''0.07476'
$ws.Range("E46").Value = '  -3.02%  '
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("E49").Value = '  +2.53%  '
$ws.Range("D50").Value = '''2.426'
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("D51").Value = '''71.93'
$ws.Range("E51").Value = '  -0.94%  '
